$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.20422986892704
$ws.Range("C2").Value = 10.21172724742361
$ws.Range("D2").Value = 4.83238250962721
$ws.Range("E2").Value = 12.21280148691869
$ws.Range("F2").Value = 24.66609776886317
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 22.11756352582325
$ws.Range("L2").Value = 10.02367379382293
$ws.Range("M2").Value = 14.38219852680218
$ws.Range("N2").Value = 17.83463344637391
$ws.Range("O2").Value = 22.00045793095433
$ws.Range("B3").Value = 13.73990935458
$ws.Range("C3").Value = 9.980055976118397
$ws.Range("D3").Value = 4.793411487157504
$ws.Range("E3").Value = 12.24915560344527
$ws.Range("F3").Value = 24.64339444581201
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 22.20910413977147
$ws.Range("L3").Value = 10.03130675015695
$ws.Range("M3").Value = 14.29024941157946
$ws.Range("N3").Value = 17.88579679501933
$ws.Range("O3").Value = 22.03916748154758
$ws.Range("B4").Value = 13.44848288234643
$ws.Range("C4").Value = 9.833686894173088
$ws.Range("D4").Value = 4.768982357822702
$ws.Range("E4").Value = 12.27270465952876
$ws.Range("F4").Value = 24.63673152213481
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 22.26998999818572
$ws.Range("L4").Value = 10.0374033365541
$ws.Range("M4").Value = 14.23559042766703
$ws.Range("N4").Value = 17.91897478227963
$ws.Range("O4").Value = 22.06847598540774
$ws.Range("B5").Value = 13.32831459669356
$ws.Range("C5").Value = 9.77305532501253
$ws.Range("D5").Value = 4.758905258728173
$ws.Range("E5").Value = 12.28261052454263
$ws.Range("F5").Value = 24.63584853381387
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 22.2959762434817
$ws.Range("L5").Value = 10.04024278590899
$ws.Range("M5").Value = 14.21378646875387
$ws.Range("N5").Value = 17.93293948704045
$ws.Range("O5").Value = 22.081809496405
$ws.Range("B6").Value = 13.30828123642056
$ws.Range("C6").Value = 9.76292957508165
$ws.Range("D6").Value = 4.757224708089232
$ws.Range("E6").Value = 12.28427409769018
$ws.Range("F6").Value = 24.63581261118838
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 22.30036213564803
$ws.Range("L6").Value = 10.04073573121928
$ws.Range("M6").Value = 14.21019483447679
$ws.Range("N6").Value = 17.93528518519687
$ws.Range("O6").Value = 22.08410736927172
$ws.Range("B7").Value = 13.44686770607901
$ws.Range("C7").Value = 9.832873113424814
$ws.Range("D7").Value = 4.768846943481596
$ws.Range("E7").Value = 12.27283699957011
$ws.Range("F7").Value = 24.63671219345723
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 22.2703357037905
$ws.Range("L7").Value = 10.03744019219962
$ws.Range("M7").Value = 14.23529444567835
$ws.Range("N7").Value = 17.91916131450795
$ws.Range("O7").Value = 22.06865018209562
$ws.Range("B8").Value = 14.04555115857307
$ws.Range("C8").Value = 10.13273245453489
$ws.Range("D8").Value = 4.819051554051287
$ws.Range("E8").Value = 12.22508208881206
$ws.Range("F8").Value = 24.65676121854692
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 22.14815432106368
$ws.Range("L8").Value = 10.02601333122902
$ws.Range("M8").Value = 14.3501321853103
$ws.Range("N8").Value = 17.85190913603159
$ws.Range("O8").Value = 22.01265323368993
$ws.Range("B9").Value = 15.16207796500661
$ws.Range("C9").Value = 10.68594802290962
$ws.Range("D9").Value = 4.913361899823479
$ws.Range("E9").Value = 12.14113898361741
$ws.Range("F9").Value = 24.75364078588575
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 21.9457738828831
$ws.Range("L9").Value = 10.01476814830285
$ws.Range("M9").Value = 14.58878220473859
$ws.Range("N9").Value = 17.73397659386121
$ws.Range("O9").Value = 21.94692975560497
$ws.Range("B10").Value = 15.93870137602903
$ws.Range("C10").Value = 11.06859139950418
$ws.Range("D10").Value = 4.979908562590854
$ws.Range("E10").Value = 12.08533221235184
$ws.Range("F10").Value = 24.859566912243
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 21.81989157081329
$ws.Range("L10").Value = 10.01327586716003
$ws.Range("M10").Value = 14.77119305268245
$ws.Range("N10").Value = 17.65577435098551
$ws.Range("O10").Value = 21.92565899360904
$ws.Range("B11").Value = 16.28098764165287
$ws.Range("C11").Value = 11.23701457216891
$ws.Range("D11").Value = 5.00954143906003
$ws.Range("E11").Value = 12.06120753507894
$ws.Range("F11").Value = 24.91519212377981
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 21.76760207721439
$ws.Range("L11").Value = 10.01405776598019
$ws.Range("M11").Value = 14.85546346271525
$ws.Range("N11").Value = 17.62201874855022
$ws.Range("O11").Value = 21.92186768430121
$ws.Range("B12").Value = 16.40891090955342
$ws.Range("C12").Value = 11.29994440079852
$ws.Range("D12").Value = 5.020666903479886
$ws.Range("E12").Value = 12.05225285282173
$ws.Range("F12").Value = 24.93731375038707
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 21.74851896836587
$ws.Range("L12").Value = 10.01456296636933
$ws.Range("M12").Value = 14.88753860184049
$ws.Range("N12").Value = 17.60949703548171
$ws.Range("O12").Value = 21.92127880817935
$ws.Range("B13").Value = 16.38143727616382
$ws.Range("C13").Value = 11.28642962240497
$ws.Range("D13").Value = 5.01827516039854
$ws.Range("E13").Value = 12.05417337305319
$ws.Range("F13").Value = 24.9325026375896
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 21.75259689559633
$ws.Range("L13").Value = 10.01444487801886
$ws.Range("M13").Value = 14.88062371322829
$ws.Range("N13").Value = 17.61218222323984
$ws.Range("O13").Value = 21.92136796829674
$ws.Range("B14").Value = 16.29154648791454
$ws.Range("C14").Value = 11.24220905000544
$ws.Range("D14").Value = 5.01045868099193
$ws.Range("E14").Value = 12.06046720815274
$ws.Range("F14").Value = 24.91699096351439
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 21.76601769919405
$ws.Range("L14").Value = 10.01409514389975
$ws.Range("M14").Value = 14.85809914748852
$ws.Range("N14").Value = 17.62098335716271
$ws.Range("O14").Value = 21.92180226460557
$ws.Range("B15").Value = 16.23626222024927
$ws.Range("C15").Value = 11.2150111474657
$ws.Range("D15").Value = 5.00565826315313
$ws.Range("E15").Value = 12.06434589058174
$ws.Range("F15").Value = 24.90762695309786
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 21.77433187672922
$ws.Range("L15").Value = 10.01390812482146
$ws.Range("M15").Value = 14.84432287560496
$ws.Range("N15").Value = 17.62640824737816
$ws.Range("O15").Value = 21.92217857017644
$ws.Range("B16").Value = 15.91610097727494
$ws.Range("C16").Value = 11.05746776956614
$ws.Range("D16").Value = 4.977958780304562
$ws.Range("E16").Value = 12.08693415522663
$ws.Range("F16").Value = 24.85608036381496
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 21.82340912635954
$ws.Range("L16").Value = 10.01325408370947
$ws.Range("M16").Value = 14.7657098942162
$ws.Range("N16").Value = 17.65801689563092
$ws.Range("O16").Value = 21.9260252458118
$ws.Range("B17").Value = 15.71679118586193
$ws.Range("C17").Value = 10.95934741867722
$ws.Range("D17").Value = 4.960799446192355
$ws.Range("E17").Value = 12.10111409285311
$ws.Range("F17").Value = 24.82635537334057
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 21.85479234938854
$ws.Range("L17").Value = 10.01322634069106
$ws.Range("M17").Value = 14.71779850351856
$ws.Range("N17").Value = 17.67787311179024
$ws.Range("O17").Value = 21.92989290939635
$ws.Range("B18").Value = 15.6011235279056
$ws.Range("C18").Value = 10.90238202693039
$ws.Range("D18").Value = 4.950869883469347
$ws.Range("E18").Value = 12.10938885315837
$ws.Range("F18").Value = 24.80995936166982
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 21.87331127762778
$ws.Range("L18").Value = 10.01334793770691
$ws.Range("M18").Value = 14.69036395884104
$ws.Range("N18").Value = 17.68946513865106
$ws.Range("O18").Value = 21.93267140182984
$ws.Range("B19").Value = 15.56178714121127
$ws.Range("C19").Value = 10.8830048147435
$ws.Range("D19").Value = 4.947497724728854
$ws.Range("E19").Value = 12.1122109757119
$ws.Range("F19").Value = 24.80452869660941
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 21.87966179710974
$ws.Range("L19").Value = 10.01341276065219
$ws.Range("M19").Value = 14.68109685017626
$ws.Range("N19").Value = 17.69341943697656
$ws.Range("O19").Value = 21.93370725979202
$ws.Range("B20").Value = 15.73811549132935
$ws.Range("C20").Value = 10.96984754069138
$ws.Range("D20").Value = 4.962632321507067
$ws.Range("E20").Value = 12.099592320235
$ws.Range("F20").Value = 24.82944717934268
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 21.85140308483655
$ws.Range("L20").Value = 10.01321506419422
$ws.Range("M20").Value = 14.72288619906716
$ws.Range("N20").Value = 17.67574166489594
$ws.Range("O20").Value = 21.92942385585272
$ws.Range("B21").Value = 16.31799636207255
$ws.Range("C21").Value = 11.25522100729419
$ws.Range("D21").Value = 5.012757202828777
$ws.Range("E21").Value = 12.05861365447867
$ws.Range("F21").Value = 24.92151852108612
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 21.76205618432832
$ws.Range("L21").Value = 10.01419220197527
$ws.Range("M21").Value = 14.86471088904236
$ws.Range("N21").Value = 17.61839117873738
$ws.Range("O21").Value = 21.92165171745318
$ws.Range("B22").Value = 16.68707065303467
$ws.Range("C22").Value = 11.43676914436558
$ws.Range("D22").Value = 5.044956198173853
$ws.Range("E22").Value = 12.03288531370404
$ws.Range("F22").Value = 24.98784998297334
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 21.70784773429585
$ws.Range("L22").Value = 10.01604918593994
$ws.Range("M22").Value = 14.95834581128235
$ws.Range("N22").Value = 17.58242912699107
$ws.Range("O22").Value = 21.92150787673165
$ws.Range("B23").Value = 16.49102901586567
$ws.Range("C23").Value = 11.34033871521429
$ws.Range("D23").Value = 5.027823542241041
$ws.Range("E23").Value = 12.04652083259587
$ws.Range("F23").Value = 24.95188865502252
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 21.73639609841308
$ws.Range("L23").Value = 10.01494693143508
$ws.Range("M23").Value = 14.90829195270227
$ws.Range("N23").Value = 17.60148393666327
$ws.Range("O23").Value = 21.92113299575955
$ws.Range("B24").Value = 15.72847813384282
$ws.Range("C24").Value = 10.96510216049138
$ws.Range("D24").Value = 4.96180387897139
$ws.Range("E24").Value = 12.10027993193449
$ws.Range("F24").Value = 24.82804721341603
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 21.85293388787325
$ws.Range("L24").Value = 10.01321973378572
$ws.Range("M24").Value = 14.72058570650349
$ws.Range("N24").Value = 17.67670474241867
$ws.Range("O24").Value = 21.92963418645116
$ws.Range("B25").Value = 14.86713266287902
$ws.Range("C25").Value = 10.54029433400187
$ws.Range("D25").Value = 4.88831454383608
$ws.Range("E25").Value = 12.16281408709211
$ws.Range("F25").Value = 24.72129841209234
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 21.99652661952783
$ws.Range("L25").Value = 10.01661855879593
$ws.Range("M25").Value = 14.5228951358726
$ws.Range("N25").Value = 17.76439351595167
$ws.Range("O25").Value = 21.92113299575955
